# Actualización automática 2025-07-31 12:55:09
#
# Updates the sales figures for asesor "GUERRERO FAREZ FABIAN MAURICIO":
#   - RUIZ TINIZARAY YOHANNA MARYURI now has a 240X120 PORCELANATO sale of 1021.25
#   - TORO BLACIO MARIA DEL CISNE now has a GRANITO sale of 153.5
# These feed through the monthly-sales summary and the compliance sheet totals.

$wb = $excel.ActiveWorkbook

$wsGrupo  = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")
$wsCumpl  = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# --- Sheet: VENTAS POR GRUPO -------------------------------------------------
$wsGrupo.Range("C49").Value2 = 1021.25
$wsGrupo.Range("F53").Value2 = 153.5

# Row 56 holds "<n> de 54" completion counters per product column.
$wsGrupo.Range("C56").Value2 = "5 de 54"
$wsGrupo.Range("F56").Value2 = "2 de 54"

# --- Sheet: VENTA MENSUAL -----------------------------------------------------
$wsMensual.Range("F49").Value2 = 1021.25
$wsMensual.Range("F53").Value2 = 153.5
$wsMensual.Range("F56").Value2 = 82305

# --- Sheet: CUMPLIMIENTO MENSUAL ----------------------------------------------
# Row 2 - 240X120 PORCELANATO
$wsCumpl.Range("D2").Value2 = 6702.91
$wsCumpl.Range("E2").Value2 = 3267.43304517915
$wsCumpl.Range("F2").Value2 = 0.6722847919702205

# Row 5 - GRANITO
$wsCumpl.Range("D5").Value2 = 307
$wsCumpl.Range("E5").Value2 = -68.68000000000001
$wsCumpl.Range("F5").Value2 = 1.288183954347096

# Row 19 - TOTAL
$wsCumpl.Range("D19").Value2 = 82305.00000000001
$wsCumpl.Range("E19").Value2 = 31401.45064517915
$wsCumpl.Range("F19").Value2 = 0.7238375618357192
